$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "P4TSSOP8_505_32X6_1M_A_F1_01"
$ws.Range("J226").Select()
